# sum 23 week 10 inputs
# Append new matchup rows to the "Nine" sheet, directly below the existing
# data (which currently ends at row 1046).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# New rows of data (Player_1, Points_1, Player_2, Points_2) to append.
$newRows = @(
    @(6, 6, 4, 14),
    @(4, 4, 7, 16),
    @(7, 5, 5, 15),
    @(6, 4, 5, 16),
    @(1, 14, 4, 6),
    @(4, 7, 2, 13),
    @(6, 8, 4, 12),
    @(3, 3, 4, 17),
    @(3, 13, 5, 7),
    @(3, 15, 4, 5),
    @(5, 17, 4, 3),
    @(3, 5, 5, 15)
)

# Find the first empty row right after the existing data.
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$lastRow = $r - 1

# Update the view so the newly added region is visible, matching how Excel
# scrolls down after data entry at the bottom of a long sheet.
$ws.Activate()
$topLeftRow = [Math]::Max(1, $lastRow - 24)
$excel.ActiveWindow.ScrollRow = $topLeftRow
$ws.Cells.Item($lastRow + 1, 1).Select()
